{"js": "// Merge the multi-run opening paragraph (\"Vous allez participer\u2026\") into a\n// single plain run, updating \"la constellation Pers\u00e9e\" -> \"la Constellation\n// de Bootes\" as part of the merged text (see commit: \"Print the first\n// paragraph with the name of the constellation.\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (\n    p.text &&\n    p.text.indexOf(\"Vous \") !== -1 &&\n    p.text.indexOf(\"campagne mondiale\") !== -1\n  ) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst newText =\n  \"Vous allez participer \u00e0 une campagne mondiale d\\u2019observation pour \" +\n  \"d\u00e9tecter les plus faibles \u00e9toiles visibles afin de mesurer la \" +\n  \"pollution lumineuse sur un site donn\u00e9. Partout dans le monde, en \" +\n  \"localisant et en observant la Constellation de Bootes dans le ciel \" +\n  \"nocturne et en la comparant aux cartes stellaires, les participants, \" +\n  \"apprendront comment l\\u2019\u00e9clairage, dans leur environnement local, \" +\n  \"influence la pollution lumineuse. Vos contributions \u00e0 la base de \" +\n  \"donn\u00e9es en ligne permettront de mesurer la qualit\u00e9 du ciel nocturne.\";\n\nfunction escapeXml(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\n// Read back the paragraph's own OOXML so the <w:p> opening tag (paraId,\n// rsid*, \u2026) and its <w:pPr> (style/spacing/indent/justification) are kept\n// byte-for-byte, then swap out only the run content for one bare run that\n// carries no direct character formatting \u2014 matching the target edit, which\n// collapses 28 runs (each with its own <w:rPr>) down to a single\n// <w:r><w:t>\u2026</w:t></w:r>.\nconst range = target.getRange();\nconst ooxmlResult = range.getOoxml();\nawait context.sync();\n\nconst existing = ooxmlResult.value;\nconst bodyTail = existing.substring(existing.indexOf(\"<w:body>\"));\nconst pOpenMatch = /<w:p\\b[^>]*>/.exec(bodyTail);\nconst pPrMatch = /<w:pPr>[\\s\\S]*?<\\/w:pPr>/.exec(bodyTail);\n\nif (!pOpenMatch) {\n  throw new Error(\"Could not locate paragraph opening tag in OOXML\");\n}\n\nconst pOpenTag = pOpenMatch[0];\nconst pPrXml = pPrMatch ? pPrMatch[0] : \"\";\n\nconst newParagraphXml =\n  pOpenTag + pPrXml + \"<w:r><w:t>\" + escapeXml(newText) + \"</w:t></w:r></w:p>\";\n\nconst newOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=' +\n  '\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphXml +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nrange.insertOoxml(newOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Merge the multi-run opening paragraph (\"Vous allez participer...\") into a\n# single plain run, updating \"la constellation Persee\" -> \"la Constellation\n# de Bootes\" as part of the merged text (see commit: \"Print the first\n# paragraph with the name of the constellation.\").\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Vous *\" -and $t -like \"*campagne mondiale*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Target paragraph not found\"\n}\n\n$newText = \"Vous allez participer \u00e0 une campagne mondiale d\u2019observation pour d\u00e9tecter les plus faibles \u00e9toiles visibles afin de mesurer la pollution lumineuse sur un site donn\u00e9. Partout dans le monde, en localisant et en observant la Constellation de Bootes dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l\u2019\u00e9clairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions \u00e0 la base de donn\u00e9es en ligne permettront de mesurer la qualit\u00e9 du ciel nocturne.\"\n\nfunction Escape-Xml([string]$s) {\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\n# Read back the paragraph's own WordOpenXML so the <w:p> opening tag\n# (paraId, rsid*, ...) and its <w:pPr> (style/spacing/indent/justification)\n# are kept byte-for-byte, then swap out only the run content for one bare\n# run that carries no direct character formatting -- matching the target\n# edit, which collapses 28 runs (each with its own <w:rPr>) down to a\n# single <w:r><w:t>...</w:t></w:r>.\n$rng = $target.Range\n$existing = $rng.WordOpenXML\n$bodyTail = $existing.Substring($existing.IndexOf(\"<w:body>\"))\n\n$null = $bodyTail -match '<w:p\\b[^>]*>'\n$pOpenTag = $matches[0]\n\n$pPrXml = \"\"\nif ($bodyTail -match '<w:pPr>[\\s\\S]*?</w:pPr>') {\n    $pPrXml = $matches[0]\n}\n\n$newParagraphXml = $pOpenTag + $pPrXml + \"<w:r><w:t>\" + (Escape-Xml $newText) + \"</w:t></w:r></w:p>\"\n\n$newOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=' + `\n    '\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + $newParagraphXml + '</w:body></w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$rng.InsertXML($newOoxml)\n"}
